$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "72.005.05"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +3.51%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "4.059.46"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +3.48%  "

$ws.Range("E4").Value = "  -0.02%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "523.67"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.23%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "149.54"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +3.30%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.627"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.14%  "

$ws.Range("E8").Value = "  +0.16%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.740"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +2.53%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.177"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +3.34%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0000342"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +2.50%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "46.66"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +10.27%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "10.75"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +4.46%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.690.22"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.041.54"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +2.67%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "21.51"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +9.14%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "14.35"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.07%  "

$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("E19").Value = "  -1.62%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "72.041.12"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +3.65%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "442.99"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +2.89%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "3.53"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +5.26%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "95.66"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +8.89%  "

$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "12.45"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +6.57%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "14.32"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "4.07"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -3.43%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "11.24"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +5.06%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "37.36"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +2.54%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "5.78"
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "13.62"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +3.32%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "700.28"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.45%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.130"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +2.76%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.93"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.59%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.89"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +15.81%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "67.59"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.54%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0917"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +10.25%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.446"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "41.07"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +2.95%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.61"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +22.14%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.154"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +4.18%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("E42").Value = "  -0.35%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.0490"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +2.62%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.13"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.38%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.81"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +1.21%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.54"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +4.36%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.147"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.84%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "3.21"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.38%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.000280"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +24.78%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "9.23"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +7.96%  "

$ws.Range("E51").Value = "  +1.78%  "
